# Starting to work on the introduction
#
# The opening paragraph of the Introduction is split in two:
#   - the first part is substantially rewritten/expanded, and
#   - the remainder ("Following robust and pertinent analysis...") is left
#     untouched but moved into its own paragraph, with a new blank
#     paragraph inserted ahead of it.
#
# Everything else in the document is unaffected.

$d = $word.ActiveDocument

$oldText = "In the modern age, climate change is an important phenomenon that has emerged and captivated the attention of the scientific community. Principally, the interest is with whether humans activity has contributed to significant and meaningful change to the climate system, especially since the industrial revolution. "

$newText = "In our modern age, climate change has emerged as an important phenomenon which has captivated the attention of the scientific community. Principally, the interest is with whether activity of mankind has contributed to significant and meaningful change to the climate system. With respect to this, the industrial revolution, and its associated timeline, is seen as a notable point of interest due to the widespread shift towards large scale operations, the adoption of machinery and the utilisation of fossil fuels across economies. Therefore, it serves as a useful demarcation point to observe the impacts of humans from the accelerated uptake of fossil fuels.^p^p"

$found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
Write-Output "Replaced intro paragraph: $found"
